# parameterize homeUrl for each environments.
$wb = $excel.ActiveWorkbook

$wsScenarios = $wb.Worksheets.Item("Scenarios")
$wsParameters = $wb.Worksheets.Item("Parameters")

# --- Parameters sheet: replace the hard-coded login URL with a bare
# base-url host that can be combined with each environment, and drop the
# hyperlink that pointed at the old literal URL.
$wsParameters.Cells.Item(7, 1).Value = "baseUrl"
$wsParameters.Cells.Item(7, 2).Value = "patient.qa.heal.com"
$wsParameters.Hyperlinks.Delete()

# Restore plain text styling to the cell that used to carry the
# "Hyperlink" cell style (now that it is a normal parameter row).
$wsParameters.Cells.Item(7, 2).Style = "Normal"
$wsParameters.Cells.Item(7, 2).NumberFormat = "@"

# --- Scenarios sheet: disable the test rows that depended on the old
# hard-coded URL/environment wiring until they're updated to use the
# parameterized baseUrl.
$wsScenarios.Cells.Item(4, 1).Value = "N"
$wsScenarios.Cells.Item(5, 1).Value = "N"
$wsScenarios.Cells.Item(6, 1).Value = "N"

# --- Cosmetic view-state updates (selection) to match the saved state.
[void]$wsScenarios.Range("A3:A7").Select()
[void]$wsParameters.Range("A7").Select()
